$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 831.1667
$ws.Range("I98").Value = 812.1852
$ws.Range("J98").Value = 1002
$ws.Range("K98").Value = 812.1852
$ws.Range("L98").Value = 1002
$ws.Range("M98").Value = 685.8148
$ws.Range("N98").Value = -3998

$ws.Range("H122").Value = 831.1667
$ws.Range("I122").Value = 812.1852
$ws.Range("J122").Value = 1002
$ws.Range("K122").Value = 2436.5556
$ws.Range("L122").Value = 3006
$ws.Range("M122").Value = 13.44439999999986
$ws.Range("N122").Value = -7906

$ws.Range("H138").Value = 3381.8096
$ws.Range("I138").Value = 1525.6666
$ws.Range("J138").Value = 4413
$ws.Range("K138").Value = 4576.9998
$ws.Range("L138").Value = 13239
$ws.Range("M138").Value = 563.0002000000004
$ws.Range("N138").Value = -23519


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5277.136
$ws.Range("I74").Value = 5375.476
$ws.Range("J74").Value = 5187.3477
$ws.Range("K74").Value = 5375.476
$ws.Range("L74").Value = 5187.3477
$ws.Range("M74").Value = -4501.476
$ws.Range("N74").Value = -6935.3477

$ws.Range("H77").Value = 5277.136
$ws.Range("I77").Value = 5375.476
$ws.Range("J77").Value = 5187.3477
$ws.Range("K77").Value = 26877.38
$ws.Range("L77").Value = 25936.7385
$ws.Range("M77").Value = -22509.38
$ws.Range("N77").Value = -34672.7385

$ws.Range("H88").Value = 1937.2
$ws.Range("I88").Value = 2068.6667
$ws.Range("J88").Value = 1740
$ws.Range("K88").Value = 2068.6667
$ws.Range("L88").Value = 1740
$ws.Range("M88").Value = -1662.6667
$ws.Range("N88").Value = -2552

$ws.Range("H91").Value = 1937.2
$ws.Range("I91").Value = 2068.6667
$ws.Range("J91").Value = 1740
$ws.Range("K91").Value = 2068.6667
$ws.Range("L91").Value = 1740
$ws.Range("M91").Value = -664.6667000000002
$ws.Range("N91").Value = -4548

$ws.Range("H113").Value = 31678.223
$ws.Range("J113").Value = 31678.223
$ws.Range("L113").Value = 31678.223
$ws.Range("N113").Value = -40356.223

$ws.Range("H132").Value = 4631.3
$ws.Range("I132").Value = 3544.1272
$ws.Range("J132").Value = 8617.6
$ws.Range("K132").Value = 10632.3816
$ws.Range("L132").Value = 25852.8
$ws.Range("M132").Value = -8102.381600000001
$ws.Range("N132").Value = -30912.8


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 108
$ws.Range("I22").Value = 110
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 110
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 63
$ws.Range("N22").Value = -446

$ws.Range("H107").Value = 999.6957
$ws.Range("I107").Value = 697.93335
$ws.Range("J107").Value = 1565.5
$ws.Range("K107").Value = 697.93335
$ws.Range("L107").Value = 1565.5
$ws.Range("M107").Value = 1222.06665
$ws.Range("N107").Value = -5405.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12823980
$ws.Range("I31").Value = 38462580
$ws.Range("J31").Value = 4680.769
$ws.Range("K31").Value = 38462580
$ws.Range("L31").Value = 4680.769
$ws.Range("M31").Value = -38462285
$ws.Range("N31").Value = -5270.769

$ws.Range("H34").Value = 12823980
$ws.Range("I34").Value = 38462580
$ws.Range("J34").Value = 4680.769
$ws.Range("K34").Value = 38462580
$ws.Range("L34").Value = 4680.769
$ws.Range("M34").Value = -38462378
$ws.Range("N34").Value = -5084.769

$ws.Range("H58").Value = 1341.1842
$ws.Range("I58").Value = 1347.7916
$ws.Range("J58").Value = 1329.8572
$ws.Range("K58").Value = 1347.7916
$ws.Range("L58").Value = 1329.8572
$ws.Range("M58").Value = -1144.7916
$ws.Range("N58").Value = -1735.8572

$ws.Range("H132").Value = 29417410
$ws.Range("I132").Value = 55563540
$ws.Range("J132").Value = 3015.875
$ws.Range("K132").Value = 166690620
$ws.Range("L132").Value = 9047.625
$ws.Range("M132").Value = -166688090
$ws.Range("N132").Value = -14107.625

$ws.Range("H134").Value = 1816.279
$ws.Range("I134").Value = 1448.8
$ws.Range("J134").Value = 3424
$ws.Range("K134").Value = 4346.4
$ws.Range("L134").Value = 10272
$ws.Range("M134").Value = -1811.4
$ws.Range("N134").Value = -15342

$ws.Range("H136").Value = 1341.1842
$ws.Range("I136").Value = 1347.7916
$ws.Range("J136").Value = 1329.8572
$ws.Range("K136").Value = 4043.3748
$ws.Range("L136").Value = 3989.5716
$ws.Range("M136").Value = -1493.3748
$ws.Range("N136").Value = -9089.5716


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 166667170
$ws.Range("I44").Value = 1000
$ws.Range("J44").Value = 333333340
$ws.Range("K44").Value = 3000
$ws.Range("L44").Value = 1000000020
$ws.Range("M44").Value = -2602
$ws.Range("N44").Value = -1000000816

$ws.Range("H47").Value = 476.625
$ws.Range("I47").Value = 53.25
$ws.Range("J47").Value = 900
$ws.Range("K47").Value = 159.75
$ws.Range("L47").Value = 2700
$ws.Range("M47").Value = 271.25
$ws.Range("N47").Value = -3562

$ws.Range("H113").Value = 581880.1
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 894930.94
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2684792.82
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -2689132.82


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 1825.1666
$ws.Range("J92").Value = 1825.1666
$ws.Range("L92").Value = 1825.1666
$ws.Range("N92").Value = -5569.1666

$ws.Range("H98").Value = 48000
$ws.Range("J98").Value = 48000
$ws.Range("L98").Value = 48000
$ws.Range("N98").Value = -53990

$ws.Range("H102").Value = 1986.24
$ws.Range("I102").Value = 2068.9333
$ws.Range("J102").Value = 1862.2
$ws.Range("K102").Value = 2068.9333
$ws.Range("L102").Value = 1862.2
$ws.Range("M102").Value = -446.9333000000001
$ws.Range("N102").Value = -5106.2


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2604.243
$ws.Range("I136").Value = 953.95654
$ws.Range("J136").Value = 5767.2915
$ws.Range("K136").Value = 2861.86962
$ws.Range("L136").Value = 17301.8745
$ws.Range("M136").Value = -311.8696199999999
$ws.Range("N136").Value = -22401.8745


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 39522.8
$ws.Range("J105").Value = 39522.8
$ws.Range("L105").Value = 39522.8
$ws.Range("N105").Value = -46510.8

$ws.Range("H126").Value = 902.85
$ws.Range("I126").Value = 712.25
$ws.Range("J126").Value = 1188.75
$ws.Range("K126").Value = 2136.75
$ws.Range("L126").Value = 3566.25
$ws.Range("M126").Value = 333.25
$ws.Range("N126").Value = -8506.25

$ws.Range("H132").Value = 7356215.5
$ws.Range("I132").Value = 10873251
$ws.Range("J132").Value = 2413.9092
$ws.Range("K132").Value = 32619753
$ws.Range("L132").Value = 7241.7276
$ws.Range("M132").Value = -32617223
$ws.Range("N132").Value = -12301.7276

